$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Rectilinear Desc fix ("Considere" -> "Consider") ---
$ws.Range("H3").Value2 = "Parallel lines spaced according to infill density. Each layer is printed perpendicular to the previous, resulting in low vertical bonding. Consider using new [Zig Zag](#zig-zag) infill instead."

# --- Row 5: 2D Lattice Desc rewrite + row height 60 -> 195 ---
$ws.Range("H5").Value2 = "Low-strength pattern with good flexibility. You can adjust **Angle 1** and **Angle 2** to optimize the infill for your specific model. Each angle adjusts the plane of each layer generated by the pattern. 0° is vertical."
$ws.Rows.Item(5).RowHeight = 195

# --- Row 10: Gyroid Desc addition + row height 150 -> 360 ---
$ws.Range("H10").Value2 = "Mathematical, isotropic surface providing equal strength in all directions. Excellent for strong, flexible prints and resin filling due to its interconnected structure. This pattern may require more time to slice because of all the points needed to generate each curve. If your model has complex geometry, consider using a simpler infill pattern like [TPMS-D](#tpms-d) or [Cross Hatch](#cross-hatch)."
$ws.Rows.Item(10).RowHeight = 360

# --- Row 11: TPMS-D Desc addition + row height 195 -> 255 ---
$ws.Range("H11").Value2 = "Triply Periodic Minimal Surface - D. Hybrid between [Cross Hatch](#cross-hatch) and [Gyroid](#gyroid), combining rigidity and smooth transitions. Isotropic and strong in all directions. This geometry is faster to slice than Gyroid, but slower than Cross Hatch."
$ws.Rows.Item(11).RowHeight = 255

# --- Row 15: 2D Honeycomb Desc fix + row height 300 -> 315 ---
$ws.Range("H15").Value2 = "Vertical Honeycomb pattern. Acceptable torsional stiffness. Developed for low densities structures like wings. Improve over [2D Lattice](#2d-lattice) offers same performance with lower densities.This infill includes a Overhang angle parameter to improve the point of contact between layers and reduce the risk of delamination."
$ws.Rows.Item(15).RowHeight = 315

# --- Row 16: 3D Honeycomb Desc typo fixes (mantaining/mantian -> maintaining/maintain) ---
$ws.Range("H16").Value2 = "This infill tries to generate a printable honeycomb structure by printing squares and octagons maintaining a vertical angle high enough to maintain contact with the previous layer."

# --- Row 22: Cross Hatch Desc addition + row height 90 -> 180 ---
$ws.Range("H22").Value2 = "Similar to [Gyroid](#gyroid) but with linear patterns, creating weak points at internal corners.`nEasier to slice but consider using [TPMS-D](#tpms-d) or [Gyroid](#gyroid) for better strength and flexibility."
$ws.Rows.Item(22).RowHeight = 180

# --- Row 24: Zig Zag Desc fix ("infil" -> "infill") ---
$ws.Range("H24").Value2 = "Similar to [rectilinear](#rectilinear) with consistent pattern between layers. Allows you to add a Symmetric infill Y axis for models with two symmetric parts."

# --- Row 25: "Coss Zag" -> "Cross Zag" (name + desc typo fixes) ---
$ws.Range("G25").Value2 = "Cross Zag"
$ws.Range("H25").Value2 = "Similar to [Zig Zag](#zig-zag) but displacing each layer with Infill shift step parameter."

# --- Row 26: Locked Zag Desc fix ("Adaptative" -> "Adaptive") ---
$ws.Range("H26").Value2 = "Adaptive version of [Zig Zag](#zig-zag) adding an external skin texture to interlock layers and a low material skeleton."

# --- Recalculate the workbook so the Infill table's calculated columns
#     (image/Pattern/MD, which derive from Infill+Desc) pick up the edits ---
$excel.CalculateFull()

# --- Sheet view: scroll down to show the newly-edited rows, matching the
#     author's final selection/viewport when they saved the file ---
$ws.Range("H27").Select()

